$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '246.74'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.32'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.093'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05612'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.500'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8123'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8432'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1343'
$ws.Range("E9").Value = '8WazirXWRX'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.03139'
$ws.Range("E10").Value = '9LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.02821'
$ws.Range("E11").Value = '10BitrueCoinBTR'
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09394'
$ws.Range("E12").Value = '11BitMartTokenBMX'
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001511'
$ws.Range("E13").Value = '12BitForexTokenBF'
$ws.Range("B14").Value = 'TigerCash'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.006113'
$ws.Range("E14").Value = '13TigerCashTCH'
$ws.Range("B15").Value = 'LEO'
$ws.Range("C15").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.568'
$ws.Range("E15").Value = '14LEOLEO'
$ws.Range("B16").Value = 'GateToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.054'
$ws.Range("E16").Value = '15GateTokenGT'
$ws.Range("B17").Value = 'BTSEToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.118'
$ws.Range("E17").Value = '16BTSETokenBTSE'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0005994'
$ws.Range("E18").Value = '17OneONE'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3181'
$ws.Range("E19").Value = '18BitpandaEcosystemTokenBEST'
$ws.Range("B20").Value = 'MandalaExchangeToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06954'
$ws.Range("E20").Value = '19MandalaExchangeTokenMDX'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.740'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04690'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004619'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009590'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001389'
$ws.Range("E28").Value = '27UpBotsUBXTWorstin24h'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03668'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006114'
$ws.Range("E41").Value = '40KickTokenKICKBestin24h'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1058'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002498'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008670'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002098'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001998'
